$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 438, shifting existing rows 438:467 down to 439:468.
$ws.Rows.Item(438).Insert()

# Populate the newly inserted row 438 with the new weekly price record
# (same market/region/category metadata as surrounding rows).
$ws.Cells.Item(438, 1).Value = 3
$ws.Cells.Item(438, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(438, 3).Value = "Coquimbo"
$ws.Cells.Item(438, 4).Value = 44931
$ws.Cells.Item(438, 5).Value = 5
$ws.Cells.Item(438, 6).Value = 100114013
$ws.Cells.Item(438, 7).Value = "Zanahoria"
$ws.Cells.Item(438, 8).Value = "Sin especificar"
$ws.Cells.Item(438, 9).Value = "Primera"
$ws.Cells.Item(438, 10).Value = 510
$ws.Cells.Item(438, 11).Value = 11500
$ws.Cells.Item(438, 12).Value = 12000
$ws.Cells.Item(438, 13).Value = 11755
$ws.Cells.Item(438, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(438, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(438, 16).Value = 588
$ws.Cells.Item(438, 17).Value = 20
$ws.Cells.Item(438, 18).Value = "Hortaliza"
